# Add the missing "8h 14m" time entry for Weston Straw (row 7, column B)
# and leave that cell selected, matching the author's final cursor position.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = "8h 14m"
$ws.Range("B7").Select()
